# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that get stamped when the handback
# status report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 09b6c85c... file
# (this text is shared with de-de!H3 "Correspond Handoff Datetime" below)
$wsOverview.Range("G3").Value = "2016-08-20 04:48:55"

# zh-cn sheet: 09b6c85c... row - Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-20 04:48:51"
$wsZhCn.Range("K3").Value = "2016-08-20 04:49:12"

# de-de sheet: 09b6c85c... row - Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-08-20 04:48:55"
$wsDeDe.Range("K3").Value = "2016-08-20 04:49:19"
